# Update deliverables list on the "Документи" sheet (sheet 1) and
# adjust sheet selection / active-tab state to match the author's save.

$wb = $excel.ActiveWorkbook

$wsDocs  = $wb.Worksheets.Item(1)   # "Документи"
$wsRoles = $wb.Worksheets.Item(2)   # "Роли"

# --- Cell content updates on the deliverables sheet ---------------------

# Row 20 ("Главен план за тестване"): E column now marked "e" (elaboration)
$wsDocs.Range("E20").Value = "e"

# Row 23 ("План за внедряване"): C1/C2/C3 move from "s" (stable) to
# "e" (elaboration); C4 remains "s".
$wsDocs.Range("G23").Value = "e"
$wsDocs.Range("H23").Value = "e"
$wsDocs.Range("I23").Value = "e"

# Row 24 ("Материали за инсталиране, администриране, поддръжка"): same shift.
$wsDocs.Range("G24").Value = "e"
$wsDocs.Range("H24").Value = "e"
$wsDocs.Range("I24").Value = "e"

# Row 25 ("План за обучение"): C1 "s"->"e", C2/C3 "c"->"e", C4 "c"->"s".
$wsDocs.Range("G25").Value = "e"
$wsDocs.Range("H25").Value = "e"
$wsDocs.Range("I25").Value = "e"
$wsDocs.Range("J25").Value = "s"

# Row 26 ("Материали за обучение"): C1/C2/C3 "c"->"e", C4 "c"->"s".
$wsDocs.Range("G26").Value = "e"
$wsDocs.Range("H26").Value = "e"
$wsDocs.Range("I26").Value = "e"
$wsDocs.Range("J26").Value = "s"

# --- Selection / active sheet bookkeeping --------------------------------

# Originally "Роли" (sheet 2) was the active tab with G13 selected; the
# saved workbook instead has "Документи" (sheet 1) active with J23
# selected, and "Роли" left with D13 selected.
$wsRoles.Activate()
$wsRoles.Range("D13").Select()

$wsDocs.Activate()
$wsDocs.Range("J23").Select()

Write-Output "applied deliverables/roles update"
